$wb = $excel.ActiveWorkbook

# Add a new worksheet named "za" at the end of the workbook (after the last sheet)
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "za"

# Populate the header row with the new course/schedule entry
$ws.Range("A1").Value = "COMP 220  A"
$ws.Range("B1").Value = "COMP PROGRAMMIN"
$ws.Range("C1").Value = "COMP PROGRAMMING II"
$ws.Range("D1").Value = "9:00:00"
$ws.Range("E1").Value = "9:50:00"
$ws.Range("F1").Value = "MWF"
$ws.Range("G1").Value = "STEM"

# "326" looks numeric, so a plain .Value assignment would store it as a
# number. Build it as text via a formula in a scratch cell, then
# copy/paste-values it across so the destination keeps the text type
# without touching any cell formatting/styles.
$ws.Range("Z1").Formula = '="326"'
$ws.Range("Z1").Copy()
$ws.Range("H1").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
